# Knights of the Force Academy - language.xlsx
# goblin sprites, some tweaks in act 1-1 lesson
#
# Rewrites the "newton_first_law_dlg_*" / "newton_first_law_inertia_dlg_*"
# rows (A27:B42) so the Key/Value rows line up 1:1 again (they had drifted
# out of sync with the voice-over script), adds two new force key/value
# pairs (forceFriction / forcePush), and reworks the inertia dialogue text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A27:B42 — realign keys/values and refresh copy -----------------------

$ws.Range("A27").Value = "forceFriction"
$ws.Range("B27").Value = "Force Friction"

$ws.Range("A28").Value = "forcePush"
$ws.Range("B28").Value = "Force Push"

$ws.Range("A29").Value = "forceZeroNetForce"
$ws.Range("B29").Value = "0 Net Force"

$ws.Range("A30").Value = "newton_first_law_dlg_1"
$ws.Range("B30").Value = "A net external force is the sum of all forces acting on the object. If the net external force is zero, then there is no acceleration."

$ws.Range("A31").Value = "newton_first_law_dlg_1_2"
$ws.Range("B31").Value = "I will now drop a large object from the sky to demonstrate."

$ws.Range("A32").Value = "newton_first_law_dlg_2"
$ws.Range("B32").Value = "There are two forces currently acting on this object: gravity and normal force."

$ws.Range("A33").Value = "newton_first_law_dlg_3"
$ws.Range("B33").Value = "The gravity force is pulling the object towards the Earth, which allows it to fall in the first place."
$ws.Range("B33").ClearFormats()

$ws.Range("A34").Value = "newton_first_law_dlg_3_2"
$ws.Range("B34").Value = "The normal force is pushing the object outwards from the surface."

$ws.Range("A35").Value = "newton_first_law_dlg_4"
$ws.Range("B35").Value = "Adding both the gravity and normal force will result in zero net force. Therefore the object is at rest."
$ws.Range("B35").VerticalAlignment = -4108
$ws.Range("B35").Font.Color = 0

$ws.Range("A36").Value = "newton_first_law_dlg_5"
$ws.Range("B36").Value = "It looks like the princess has been taken hostage by these nefarious goblins!"

$ws.Range("A37").Value = "newton_first_law_dlg_6"
$ws.Range("B37").Value = "Gather your forces for a daring rescue!"

$ws.Range("A38").Value = "newton_first_law_inertia_dlg_1"
$ws.Range("B38").Value = "Excellent! Notice how it took a couple of knights to push the block?"
$ws.Range("B38").ClearFormats()

$ws.Range("A39").Value = "newton_first_law_inertia_dlg_2"
$ws.Range("B39").Value = "The block's mass gives it more weight, which allows its frictional force to go against the pushing force of the knights."

$ws.Range("A40").Value = "newton_first_law_inertia_dlg_3"
$ws.Range("B40").Value = "The mass of an object determines the object’s inertia. The unit of measurement shown here is in kg (kilograms), which is 1000 grams per 1 kilogram."
$ws.Range("B40").ClearFormats()

$ws.Range("A41").Value = "newton_first_law_inertia_dlg_4"
$ws.Range("B41").Value = "Inertia is the tendency of any physical object to resist any change in motion. A change in motion is acceleration due to net force acting on the object."
$ws.Range("B41").VerticalAlignment = -4108

$ws.Range("A42").Value = "newton_first_law_inertia_dlg_5"
$ws.Range("B42").Value = "I will now drop another object from the sky. This time with less mass. Just one more block to rescue our damsel in distress!"

# --- selection moved from B41 to B39 --------------------------------------

[void]$ws.Range("B39").Select()
